$p = $ppt.ActivePresentation
$cr = [char]13

# ---------------------------------------------------------------------------
# 1. Slide 12 ("We Need More Control"): merge the "WCF " + "serializers"
#    runs into a single run "WCF serializers".
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$contentShape12 = $s12.Shapes.Item(2)
$tr12 = $contentShape12.TextFrame.TextRange
$full12 = $tr12.Text
$idx12 = $full12.IndexOf("WCF serializers")
$sub12 = $tr12.Characters($idx12 + 1, 15)
$sub12.Text = "WCF serializers"

# ---------------------------------------------------------------------------
# 2. Append four new "Title and Content" slides at the end of the deck.
# ---------------------------------------------------------------------------

# --- New slide 15: "Versioning" ---
$s15 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s15.Shapes.Item(1).TextFrame.TextRange.Text = "Versioning"

$body15 = $s15.Shapes.Item(2)
$body15.Left = 18
$body15.Top = 126
$body15.Width = 690
$body15.Height = 357

$tr15 = $body15.TextFrame.TextRange
$tr15.Text = ".NET 2.0 added version tolerance to BF" + $cr + `
    "Add new fields without breaking deserialization" + $cr + `
    "Additional fields in stream ignored during deserialization" + $cr + `
    "Attribute fields as OptionalField" + $cr + `
    "Set VersionAdded param, though still not implemented"

$tr15.Paragraphs(2, 1).IndentLevel = 2
$tr15.Paragraphs(3, 1).IndentLevel = 2
$tr15.Paragraphs(4, 1).IndentLevel = 2
$tr15.Paragraphs(5, 1).IndentLevel = 3

# --- New slide 16: "Version Tolerance Rules" ---
$s16 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s16.Shapes.Item(1).TextFrame.TextRange.Text = "Version Tolerance Rules"

$body16 = $s16.Shapes.Item(2)
$tr16 = $body16.TextFrame.TextRange
$tr16.Text = "Never remove a serialized field." + $cr + `
    "Never apply the NonSerializedAttribute attribute to a field if the attribute was not applied to the field in the previous version." + $cr + `
    "Never change the name or the type of a serialized field." + $cr + `
    "When adding a new serialized field, apply the OptionalFieldAttribute attribute."

# --- New slide 17: "Version Tolerance Rules" (continued) ---
$s17 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s17.Shapes.Item(1).TextFrame.TextRange.Text = "Version Tolerance Rules"

$body17 = $s17.Shapes.Item(2)
$tr17 = $body17.TextFrame.TextRange
$tr17.Text = "When removing a NonSerializedAttribute attribute from a field (that was not serializable in a previous version), apply the OptionalFieldAttribute attribute." + $cr + `
    "For all optional fields, set meaningful defaults using the serialization callbacks unless 0 or nullas defaults are acceptable." + $cr + `
    "x"

# Trim the placeholder character on the trailing empty paragraph so it is a
# genuinely empty paragraph (no stray run) instead of a run containing "x".
$full17 = $tr17.Text
$lastChar17 = $tr17.Characters($full17.Length, 1)
$lastChar17.Text = ""

# --- New slide 18: "Versioning" (continued) ---
$s18 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s18.Shapes.Item(1).TextFrame.TextRange.Text = "Versioning"

$body18 = $s18.Shapes.Item(2)
$tr18 = $body18.TextFrame.TextRange
$tr18.Text = "Other changes and other serializers require custom serialization" + $cr + `
    "Member data type changes" + $cr + `
    "Semantic changes"

$tr18.Paragraphs(2, 1).IndentLevel = 2
$tr18.Paragraphs(3, 1).IndentLevel = 2
